# Scotland League Two - corrective re-shuffle of fixture rows.
#
# The source data for several match-days had rows in the wrong order
# (rows sharing the same Date got shuffled). This script restores the
# correct row order by rotating/swapping the row contents (columns B..AD)
# while leaving column A (the positional index) untouched.
#
# Implementation: snapshot every involved row's B..AD values first (so
# reads never see a partially-updated row), then write each destination
# row from its mapped source row's snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

function Get-RowSnapshot($r) {
    $data = @{}
    foreach ($c in $cols) {
        $data[$c] = $ws.Range("$c$r").Value2
    }
    return $data
}

function Set-RowFromSnapshot($r, $data) {
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $data[$c]
    }
}

# destination row -> source row (source row's CURRENT content is copied into destination row)
$mapping = @{
    22 = 23; 23 = 22;
    42 = 43; 43 = 44; 44 = 42;
    52 = 53; 53 = 52;
    80 = 83; 83 = 80;
    87 = 88; 88 = 89; 89 = 90; 90 = 91; 91 = 87;
    105 = 106; 106 = 105;
    108 = 111; 109 = 110; 110 = 109; 111 = 108;
    123 = 124; 124 = 123;
}

# Snapshot every row that participates (as source or destination) before any writes.
$snapshots = @{}
foreach ($r in $mapping.Keys) {
    if (-not $snapshots.ContainsKey($r)) {
        $snapshots[$r] = Get-RowSnapshot $r
    }
}

# Apply: write destination row using the snapshot of its source row.
foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]
    Set-RowFromSnapshot $dest $snapshots[$src]
}
